$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Wnt11"
$ws.Cells.Item(2,3).Value = "Fzd8"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.1470843333333333
$ws.Cells.Item(2,8).Value = 0.441253
$ws.Cells.Item(2,9).Value = 0.01080519019772543
$ws.Cells.Item(2,10).Value = 0.01080519019772543
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 2.558821666666667
$ws.Cells.Item(2,14).Value = 7.676465
$ws.Cells.Item(2,15).Value = 0.2156728774407755
$ws.Cells.Item(2,16).Value = 0.2156728774407755
$ws.Cells.Item(2,17).Value = 0.3763625789605556
$ws.Cells.Item(2,18).Value = 3.387263210645
$ws.Cells.Item(2,19).Value = 0.002330386461238306
$ws.Cells.Item(2,20).Value = 0.002330386461238306

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Wnt11"
$ws.Cells.Item(3,3).Value = "Fzd8"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.1470843333333333
$ws.Cells.Item(3,8).Value = 0.441253
$ws.Cells.Item(3,9).Value = 0.01080519019772543
$ws.Cells.Item(3,10).Value = 0.01080519019772543
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 6.453984666666667
$ws.Cells.Item(3,14).Value = 19.361954
$ws.Cells.Item(3,15).Value = 0.5439806384912759
$ws.Cells.Item(3,16).Value = 0.5439806384912759
$ws.Cells.Item(3,17).Value = 0.9492800320402224
$ws.Cells.Item(3,18).Value = 8.543520288362
$ws.Cells.Item(3,19).Value = 0.005877814262778357
$ws.Cells.Item(3,20).Value = 0.005877814262778357

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Wnt11"
$ws.Cells.Item(4,3).Value = "Fzd8"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.1470843333333333
$ws.Cells.Item(4,8).Value = 0.441253
$ws.Cells.Item(4,9).Value = 0.01080519019772543
$ws.Cells.Item(4,10).Value = 0.01080519019772543
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.851558333333334
$ws.Cells.Item(4,14).Value = 8.554675000000001
$ws.Cells.Item(4,15).Value = 0.2403464840679487
$ws.Cells.Item(4,16).Value = 0.2403464840679487
$ws.Cells.Item(4,17).Value = 0.4194195564194445
$ws.Cells.Item(4,18).Value = 3.774776007775001
$ws.Cells.Item(4,19).Value = 0.002596989473708771
$ws.Cells.Item(4,20).Value = 0.002596989473708772

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Wnt11"
$ws.Cells.Item(5,3).Value = "Fzd8"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 13.22273466666667
$ws.Cells.Item(5,8).Value = 39.668204
$ws.Cells.Item(5,9).Value = 0.9713758071269154
$ws.Cells.Item(5,10).Value = 0.9713758071269153
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 2.558821666666667
$ws.Cells.Item(5,14).Value = 7.676465
$ws.Cells.Item(5,15).Value = 0.2156728774407755
$ws.Cells.Item(5,16).Value = 0.2156728774407755
$ws.Cells.Item(5,17).Value = 33.83461995765111
$ws.Cells.Item(5,18).Value = 304.51157961886
$ws.Cells.Item(5,19).Value = 0.2094994153994176
$ws.Cells.Item(5,20).Value = 0.2094994153994175

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Wnt11"
$ws.Cells.Item(6,3).Value = "Fzd8"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 13.22273466666667
$ws.Cells.Item(6,8).Value = 39.668204
$ws.Cells.Item(6,9).Value = 0.9713758071269154
$ws.Cells.Item(6,10).Value = 0.9713758071269153
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 6.453984666666667
$ws.Cells.Item(6,14).Value = 19.361954
$ws.Cells.Item(6,15).Value = 0.5439806384912759
$ws.Cells.Item(6,16).Value = 0.5439806384912759
$ws.Cells.Item(6,17).Value = 85.33932679006844
$ws.Cells.Item(6,18).Value = 768.0539411106159
$ws.Cells.Item(6,19).Value = 0.5284096317758779
$ws.Cells.Item(6,20).Value = 0.5284096317758779

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Wnt11"
$ws.Cells.Item(7,3).Value = "Fzd8"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 13.22273466666667
$ws.Cells.Item(7,8).Value = 39.668204
$ws.Cells.Item(7,9).Value = 0.9713758071269154
$ws.Cells.Item(7,10).Value = 0.9713758071269153
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.851558333333334
$ws.Cells.Item(7,14).Value = 8.554675000000001
$ws.Cells.Item(7,15).Value = 0.2403464840679487
$ws.Cells.Item(7,16).Value = 0.2403464840679487
$ws.Cells.Item(7,17).Value = 37.70539922818889
$ws.Cells.Item(7,18).Value = 339.3485930537
$ws.Cells.Item(7,19).Value = 0.23346675995162
$ws.Cells.Item(7,20).Value = 0.23346675995162

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Wnt11"
$ws.Cells.Item(8,3).Value = "Fzd8"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.242559
$ws.Cells.Item(8,8).Value = 0.7276769999999999
$ws.Cells.Item(8,9).Value = 0.01781900267535914
$ws.Cells.Item(8,10).Value = 0.01781900267535914
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 2.558821666666667
$ws.Cells.Item(8,14).Value = 7.676465
$ws.Cells.Item(8,15).Value = 0.2156728774407755
$ws.Cells.Item(8,16).Value = 0.2156728774407755
$ws.Cells.Item(8,17).Value = 0.620665224645
$ws.Cells.Item(8,18).Value = 5.585987021805
$ws.Cells.Item(8,19).Value = 0.003843075580119583
$ws.Cells.Item(8,20).Value = 0.003843075580119583

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Wnt11"
$ws.Cells.Item(9,3).Value = "Fzd8"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.242559
$ws.Cells.Item(9,8).Value = 0.7276769999999999
$ws.Cells.Item(9,9).Value = 0.01781900267535914
$ws.Cells.Item(9,10).Value = 0.01781900267535914
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 6.453984666666667
$ws.Cells.Item(9,14).Value = 19.361954
$ws.Cells.Item(9,15).Value = 0.5439806384912759
$ws.Cells.Item(9,16).Value = 0.5439806384912759
$ws.Cells.Item(9,17).Value = 1.565472066762
$ws.Cells.Item(9,18).Value = 14.089248600858
$ws.Cells.Item(9,19).Value = 0.00969319245261962
$ws.Cells.Item(9,20).Value = 0.00969319245261962

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Wnt11"
$ws.Cells.Item(10,3).Value = "Fzd8"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.242559
$ws.Cells.Item(10,8).Value = 0.7276769999999999
$ws.Cells.Item(10,9).Value = 0.01781900267535914
$ws.Cells.Item(10,10).Value = 0.01781900267535914
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.851558333333334
$ws.Cells.Item(10,14).Value = 8.554675000000001
$ws.Cells.Item(10,15).Value = 0.2403464840679487
$ws.Cells.Item(10,16).Value = 0.2403464840679487
$ws.Cells.Item(10,17).Value = 0.691671137775
$ws.Cells.Item(10,18).Value = 6.225040239975
$ws.Cells.Item(10,19).Value = 0.004282734642619941
$ws.Cells.Item(10,20).Value = 0.004282734642619942
